$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 127

$ws.Range("E8").Value = 7

$ws.Range("E10").Value = 445
$ws.Range("F10").Value = 219
$ws.Range("H10").Value = 219

$ws.Range("E11").Value = 305
$ws.Range("F11").Value = 169
$ws.Range("H11").Value = 169

$ws.Range("E12").Value = 436

$ws.Range("E13").Value = 115
$ws.Range("F13").Value = 60
$ws.Range("H13").Value = 60

$ws.Range("E14").Value = 112
$ws.Range("F14").Value = 57
$ws.Range("H14").Value = 57

$ws.Range("E15").Value = 142

$ws.Range("E16").Value = 182

$ws.Range("E22").Value = 154

$ws.Range("E24").Value = 188
$ws.Range("F24").Value = 99
$ws.Range("H24").Value = 99

$ws.Range("E25").Value = 235

$ws.Range("E27").Value = 301
$ws.Range("F27").Value = 141
$ws.Range("H27").Value = 141

$ws.Range("F29").Value = 83
$ws.Range("H29").Value = 83

$ws.Range("E30").Value = 191
$ws.Range("F30").Value = 109
$ws.Range("H30").Value = 109

$ws.Range("E32").Value = 170

$ws.Range("E35").Value = 132
$ws.Range("F35").Value = 84
$ws.Range("H35").Value = 84

$ws.Range("E36").Value = 63

$ws.Range("E40").Value = 238

$ws.Range("E41").Value = 361
$ws.Range("F41").Value = 169
$ws.Range("H41").Value = 169

$ws.Range("E42").Value = 328
$ws.Range("F42").Value = 181
$ws.Range("H42").Value = 181

$ws.Range("E44").Value = 283
$ws.Range("F44").Value = 138
$ws.Range("H44").Value = 138

$ws.Range("E45").Value = 125
$ws.Range("F45").Value = 63
$ws.Range("H45").Value = 63

$ws.Range("E46").Value = 284
$ws.Range("F46").Value = 154
$ws.Range("H46").Value = 154

$ws.Range("E47").Value = 400
$ws.Range("F47").Value = 197
$ws.Range("H47").Value = 197

$ws.Range("E48").Value = 186
$ws.Range("F48").Value = 77
$ws.Range("H48").Value = 77

$ws.Range("E49").Value = 264
$ws.Range("F49").Value = 111
$ws.Range("H49").Value = 111

$ws.Range("E51").Value = 215
